$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'knee high tights black',
    'knee highs black',
    'knee highs stockings',
    'knee hiking support',
    'knee joint support',
    'knee joint support brace',
    'knee joints',
    'knee kinesiology tape',
    'knee leggings',
    'knee leggings for basketball',
    'knee lenght shorts',
    'knee lenght tights',
    'knee length basketball shorts for women',
    'knee length bike shorts',
    'knee length compression garment',
    'knee length down jacket women',
    'knee length leggings',
    'knee length leggings for women',
    'knee length leggings with pockets for women',
    'knee length leggings women',
    'knee length pants',
    'knee length pants women',
    'knee length shorts',
    'knee length shorts for women',
    'knee length stockings for women',
    'knee length tights',
    'knee length tights men',
    'knee length tights women',
    'knee length top',
    'knee length tops for women',
    'knee length winter jackets for women',
    'knee length workout shorts',
    'knee length workout shorts for women',
    'knee length yoga pants',
    'knee length yoga pants for women',
    'knee models',
    'knee nike',
    'knee pad compression pants',
    'knee patches for work pants',
    'knee powers joint support',
    'knee pro',
    'knee protection for running',
    'knee protection for women',
    'knee race',
    'knee recovery',
    'knee reinforcement for pants',
    'knee reinforcement patch',
    'knee replacement recovery',
    'knee replacement recovery gifts',
    'knee running',
    'knee running band',
    'knee running support',
    'knee saver',
    'knee savers',
    'knee savers red',
    'knee shorts be',
    'knee shorts women',
    'knee skins',
    'knee sleep support',
    'knee sleeve red line',
    'knee sleeves exo',
    'knee sleeves exo sleeve',
    'knee sport support',
    'knee sports tape',
    'knee straightening',
    'knee supplement for women',
    'knee support basketball',
    'knee support big size',
    'knee support black',
    'knee support blue',
    'knee support compression',
    'knee support cycling',
    'knee support extra large',
    'knee support for bed',
    'knee support for runners',
    'knee support for running',
    'knee support for running for women',
    'knee support for running men',
    'knee support for running women',
    'knee support for skiing',
    'knee support for sports',
    'knee support for weightlifting women',
    'knee support for women running',
    'knee support for work',
    'knee support gym',
    'knee support hiking',
    'knee support hinge',
    'knee support ladies',
    'knee support leggings',
    'knee support pair',
    'knee support pair for women',
    'knee support pants',
    'knee support pants men',
    'knee support patella',
    'knee support plus size',
    'knee support plus size women',
    'knee support runners',
    'knee support running',
    'knee support running women',
    'knee support skiing'
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
